$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Preserve the existing header style (cell style index used by A1:I1,
# the bold/bordered/centered look) by stashing a formatted copy in an
# unused scratch cell before we wipe the sheet.
# ------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Wipe the old table (A1:I6) entirely - content and formatting.
# ------------------------------------------------------------------
$ws.Range("A1:I6").Clear()

# ------------------------------------------------------------------
# Row 1 index cells (A-E) - plain numbers, no shared strings involved.
# ------------------------------------------------------------------
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# ------------------------------------------------------------------
# Introduce each distinct label text in the same order the source
# workbook first defines it, so the shared-string table comes out in
# the same sequence.
# ------------------------------------------------------------------
$ws.Range("F1").Value = "NI: `$\hat\sigma_{pb,SPF}`$"   # new string 0
$ws.Range("G1").Value = "`$\hat\sigma_{pr,SPF}`$"        # new string 1
$ws.Range("J1").Value = "NI: `$\rho`$"                   # new string 2
$ws.Range("K1").Value = "NI: `$\sigma`$"                 # new string 3
$ws.Range("L1").Value = "NI: `$\hat\sigma_{pb,SCE}`$"    # new string 4
$ws.Range("M1").Value = "`$\hat\sigma_{pr,SCE}`$"        # new string 5

$ws.Range("A2").Value = "FEVar"                          # new string 6
$ws.Range("A3").Value = "DisgVar"                        # new string 7
$ws.Range("B2").Value = "FEATV"                           # new string 8
$ws.Range("B3").Value = "DisgATV"                         # new string 9
$ws.Range("E5").Value = "Var"                             # new string 10

# Remaining cells reuse one of the eleven strings above.
$ws.Range("H1").Value = "NI: `$\hat\sigma_{pb,SPF}`$"
$ws.Range("I1").Value = "`$\hat\sigma_{pr,SPF}`$"
$ws.Range("N1").Value = "NI: `$\hat\sigma_{pb,SCE}`$"
$ws.Range("O1").Value = "`$\hat\sigma_{pr,SCE}`$"
$ws.Range("P1").Value = "NI: `$\rho`$"
$ws.Range("Q1").Value = "NI: `$\sigma`$"

$ws.Range("A4").Value = "FEVar"
$ws.Range("B4").Value = "FEATV"
$ws.Range("C4").Value = "DisgVar"
$ws.Range("D4").Value = "DisgATV"

$ws.Range("A5").Value = "FEVar"
$ws.Range("B5").Value = "FEATV"
$ws.Range("C5").Value = "DisgVar"
$ws.Range("D5").Value = "DisgATV"

# ------------------------------------------------------------------
# Re-apply the stashed header style to the full new header row, then
# drop the scratch cell.
# ------------------------------------------------------------------
$ws.Range("A20").Copy()
$ws.Range("A1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A20").Clear()

# ------------------------------------------------------------------
# Numeric data, row 2
# ------------------------------------------------------------------
$ws.Range("F2").Value = 25.32
$ws.Range("G2").Value = 16.07
$ws.Range("H2").Value = 0.9
$ws.Range("I2").Value = 1.09
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3.4
$ws.Range("M2").Value = 15.4
$ws.Range("N2").Value = 3.4
$ws.Range("O2").Value = 11.29
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0

# ------------------------------------------------------------------
# Numeric data, row 3
# ------------------------------------------------------------------
$ws.Range("F3").Value = 471301.73
$ws.Range("G3").Value = 0.85
$ws.Range("H3").Value = -1.74
$ws.Range("I3").Value = -0.17
$ws.Range("J3").Value = 0.91
$ws.Range("K3").Value = 0.42
$ws.Range("L3").Value = 168519.46
$ws.Range("M3").Value = 1.09
$ws.Range("N3").Value = 0.67
$ws.Range("O3").Value = 0.58
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.03

# ------------------------------------------------------------------
# Numeric data, row 4
# ------------------------------------------------------------------
$ws.Range("F4").Value = 25.32
$ws.Range("G4").Value = 16.07
$ws.Range("H4").Value = 0.9
$ws.Range("I4").Value = 1.09
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3.4
$ws.Range("M4").Value = 15.4
$ws.Range("N4").Value = 3.4
$ws.Range("O4").Value = 11.29
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0

# ------------------------------------------------------------------
# Numeric data, row 5
# ------------------------------------------------------------------
$ws.Range("F5").Value = 9167076061667.699
$ws.Range("G5").Value = 2.37
$ws.Range("H5").Value = 0.9
$ws.Range("I5").Value = 1.09
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 13.22
$ws.Range("M5").Value = 59.96
$ws.Range("N5").Value = 13.22
$ws.Range("O5").Value = 59.96
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = -0.03
